# Append: 2026-01-15 06:30 JST
# Replace the two top rows of the "ランサーズ" sheet with the newly
# scraped entries and drop the rest of the old listing (rows 4-18),
# mirroring the scraper's "latest snapshot" behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the obsolete listing rows (old rows 4 through 18) -------------
$ws.Rows("4:18").Delete()

# --- Row 2: now the former "PMO" listing (was row 7) -----------------------
$ws.Range("A2").Value = "2026-01-15 06:30:32"
$ws.Range("B2").Value = "【長期案件】生成AIを利用したチャットボット作成のPMOを募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5471035"
$ws.Range("G2").Value = 310
$ws.Range("H2").Value = "🔥AI,Ai"

# --- Row 3: brand-new listing ----------------------------------------------
$ws.Range("A3").Value = "2026-01-15 06:30:32"
$ws.Range("B3").Value = "【急募】GASプロジェクトの作成依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5471552"
$ws.Range("G3").Value = 25
$ws.Range("H3").ClearContents()

# --- Rebuild hyperlinks so only F2/F3 point at the new URLs ---------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5471035")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5471552")
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"

# --- Column width tweaks (B: 54 -> 33, D: 32 -> 30) ------------------------
# Excel's ColumnWidth property stores a little extra padding (~0.8333 chars)
# on top of whatever is assigned, so back that padding out to land on the
# exact target widths recorded in the sheet XML.
$ws.Columns.Item(2).ColumnWidth = 33 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 30 - 0.8333333333333334
